$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 52, shifting rows 52:73 down to 53:74
$ws.Rows.Item(52).Insert()

# Fill the new row 52 with data (copy of row 53's data pattern, with date/price changes)
$ws.Cells.Item(52, 1).Value = 1
$ws.Cells.Item(52, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(52, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(52, 4).Value = 44784
$ws.Cells.Item(52, 4).NumberFormat = $ws.Cells.Item(53, 4).NumberFormat
$ws.Cells.Item(52, 5).Value = 15
$ws.Cells.Item(52, 6).Value = 100112012
$ws.Cells.Item(52, 7).Value = "Espinaca"
$ws.Cells.Item(52, 8).Value = "Sin especificar"
$ws.Cells.Item(52, 9).Value = "Primera"
$ws.Cells.Item(52, 10).Value = 250
$ws.Cells.Item(52, 11).Value = 3000
$ws.Cells.Item(52, 12).Value = 3500
$ws.Cells.Item(52, 13).Value = 3250
$ws.Cells.Item(52, 14).Value = "`$/atado 2,5 a 3 kilos"
$ws.Cells.Item(52, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(52, 16).Value = 1083
$ws.Cells.Item(52, 17).Value = 3
$ws.Cells.Item(52, 18).Value = "Hortaliza"
